$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Tipo de producto" text for row 2 (A2)
$ws.Range("A2").Value = "2-Queque,1-Torta,"

# Update "Estado" (status) column G for rows 3-9 to reflect
# the current order progress (generating data for the "most sold" report chart)
$ws.Range("G3").Value = "Finalizado"
$ws.Range("G4").Value = "Finalizado"
$ws.Range("G5").Value = "Retirado"
$ws.Range("G6").Value = "Retirado"
$ws.Range("G7").Value = "Retirado"
$ws.Range("G8").Value = "Retirado"
$ws.Range("G9").Value = "Retirado"
